$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column D (NORM) updates on the ELEMENTS sheet ---
# Row 2: Tee -> now references a dedicated "Type B" norm string
$ws.Range("D2").Value = "''EN 10253-2 - Type B'"

# Rows 3-4: Reducers keep the EN 10253-2 norm, but it is now quote-prefixed text
$ws.Range("D3").Value = "''EN 10253-2'"
$ws.Range("D4").Value = "''EN 10253-2'"

# Rows 5-14: Flanges / blind flanges move from "EN 1092-1 A1" to "EN 1092-1/11/PN40"
$ws.Range("D5").Value  = "''EN 1092-1/11/PN40'"
$ws.Range("D6").Value  = "''EN 1092-1/11/PN40'"
$ws.Range("D7").Value  = "''EN 1092-1/11/PN40'"
$ws.Range("D8").Value  = "''EN 1092-1/11/PN40'"
$ws.Range("D9").Value  = "''EN 1092-1/11/PN40'"
$ws.Range("D10").Value = "''EN 1092-1/11/PN40'"
$ws.Range("D11").Value = "''EN 1092-1/11/PN40'"
$ws.Range("D12").Value = "''EN 1092-1/11/PN40'"
$ws.Range("D13").Value = "''EN 1092-1/11/PN40'"
$ws.Range("D14").Value = "''EN 1092-1/11/PN40'"

# Rows 15-17 (Elbows) keep "EN 10253-2" unchanged.

# --- Column width / view changes ---
$ws.Columns.Item(4).ColumnWidth = 18.33

$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("H18").Select()
